# Refresh of the bat-survey export rows (A2:R7) — species records were
# re-fetched/re-sorted upstream, so each row's identifiers, species data and
# location coordinates are updated in place to the new source order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Force text formatting so numeric-looking strings (e.g. the "Antal"
    # column) stay stored as text rather than being auto-coerced to numbers.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

# -- Row 2 ------------------------------------------------------------------
$ws.Range("A2").Value = 111545401
$ws.Range("B2").Value = 57487
Set-TextCell "D2" "NT"
$ws.Range("E2").Value = 205998
Set-TextCell "F2" "Nordfladdermus"
Set-TextCell "G2" "Eptesicus nilssonii"
Set-TextCell "H2" "(A.Keyserling & Blasius, 1839)"
Set-TextCell "I2" "6"
Set-TextCell "P2" "Orsa Viborg, glänta i mitten av skogsparti, Dlr"
$ws.Range("Q2").Value = 480487.2503558649
$ws.Range("R2").Value = 6772784.264016891

# -- Row 3 ------------------------------------------------------------------
$ws.Range("A3").Value = 111545414
$ws.Range("B3").Value = 57494
Set-TextCell "D3" "LC"
$ws.Range("E3").Value = 205992
Set-TextCell "F3" "Vattenfladdermus"
Set-TextCell "G3" "Myotis daubentonii"
Set-TextCell "H3" "(Kuhl, 1817)"
Set-TextCell "I3" "9"

# -- Row 4 ------------------------------------------------------------------
$ws.Range("A4").Value = 111545323
$ws.Range("B4").Value = 57487
Set-TextCell "D4" "NT"
$ws.Range("E4").Value = 205998
Set-TextCell "F4" "Nordfladdermus"
Set-TextCell "G4" "Eptesicus nilssonii"
Set-TextCell "H4" "(A.Keyserling & Blasius, 1839)"
Set-TextCell "I4" "2"

# -- Row 5 ------------------------------------------------------------------
$ws.Range("A5").Value = 111543968
$ws.Range("B5").Value = 57487
Set-TextCell "D5" "NT"
$ws.Range("E5").Value = 205998
Set-TextCell "F5" "Nordfladdermus"
Set-TextCell "G5" "Eptesicus nilssonii"
Set-TextCell "H5" "(A.Keyserling & Blasius, 1839)"
Set-TextCell "I5" "256"
Set-TextCell "J5" ""
Set-TextCell "P5" "Orsa Viborg, glänta i skogsparti, Dlr"
$ws.Range("Q5").Value = 480406.6045043401
$ws.Range("R5").Value = 6772745.04339793

# -- Row 6 ------------------------------------------------------------------
$ws.Range("A6").Value = 111545328
$ws.Range("B6").Value = 57494
Set-TextCell "D6" "LC"
$ws.Range("E6").Value = 205992
Set-TextCell "F6" "Vattenfladdermus"
Set-TextCell "G6" "Myotis daubentonii"
Set-TextCell "H6" "(Kuhl, 1817)"
Set-TextCell "I6" "1"

# -- Row 7 ------------------------------------------------------------------
$ws.Range("A7").Value = 111543957
$ws.Range("B7").Value = 57494
Set-TextCell "D7" "LC"
$ws.Range("E7").Value = 205992
Set-TextCell "F7" "Vattenfladdermus"
Set-TextCell "G7" "Myotis daubentonii"
Set-TextCell "H7" "(Kuhl, 1817)"
Set-TextCell "I7" "1"
Set-TextCell "J7" "registreringar"
